$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1289.5555
$ws.Range("I28").Value = 1300.1333
$ws.Range("J28").Value = 1236.6666
$ws.Range("K28").Value = 1300.1333
$ws.Range("L28").Value = 1236.6666
$ws.Range("M28").Value = -815.1333
$ws.Range("N28").Value = -2206.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3797.1667
$ws.Range("I76").Value = 4157.6
$ws.Range("K76").Value = 4157.6
$ws.Range("M76").Value = -3842.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3797.1667
$ws.Range("I79").Value = 4157.6
$ws.Range("K79").Value = 4157.6
$ws.Range("M79").Value = -3065.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 462.75
$ws.Range("J96").Value = 163
$ws.Range("L96").Value = 489
$ws.Range("N96").Value = -3235

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1223.2273
$ws.Range("I100").Value = 1348.5264
$ws.Range("J100").Value = 429.66666
$ws.Range("K100").Value = 1348.5264
$ws.Range("L100").Value = 429.66666
$ws.Range("M100").Value = -807.5264
$ws.Range("N100").Value = -1511.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 195.44444
$ws.Range("I101").Value = 198.625
$ws.Range("K101").Value = 595.875
$ws.Range("M101").Value = 1026.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 971.94116
$ws.Range("I107").Value = 965.3333
$ws.Range("J107").Value = 987.8
$ws.Range("K107").Value = 965.3333
$ws.Range("L107").Value = 987.8
$ws.Range("M107").Value = 954.6667
$ws.Range("N107").Value = -4827.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 14061.5
$ws.Range("J111").Value = 17391.572
$ws.Range("L111").Value = 52174.716
$ws.Range("N111").Value = -58308.716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 29021.486
$ws.Range("I112").Value = 1390.1666
$ws.Range("J112").Value = 34045.363
$ws.Range("K112").Value = 4170.4998
$ws.Range("L112").Value = 102136.089
$ws.Range("M112").Value = -3062.4998
$ws.Range("N112").Value = -104352.089

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 309.0909
$ws.Range("I115").Value = 309.0909
$ws.Range("K115").Value = 927.2727
$ws.Range("M115").Value = 639.7273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7285.5713
$ws.Range("I116").Value = 7599.8
$ws.Range("K116").Value = 7599.8
$ws.Range("M116").Value = -4157.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1848.091
$ws.Range("I137").Value = 1892.0667
$ws.Range("J137").Value = 1811.4445
$ws.Range("K137").Value = 5676.2001
$ws.Range("L137").Value = 5434.333500000001
$ws.Range("M137").Value = -3126.2001
$ws.Range("N137").Value = -10534.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1306.3
$ws.Range("I141").Value = 1217.5555
$ws.Range("J141").Value = 2105
$ws.Range("K141").Value = 3652.6665
$ws.Range("L141").Value = 6315
$ws.Range("M141").Value = 1527.3335
$ws.Range("N141").Value = -16675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1559.3334
$ws.Range("I2").Value = 1506
$ws.Range("J2").Value = 1906
$ws.Range("K2").Value = 1506
$ws.Range("L2").Value = 1906
$ws.Range("M2").Value = -1393
$ws.Range("N2").Value = -2132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 4499.5
$ws.Range("J15").Value = 4499.5
$ws.Range("L15").Value = 4499.5
$ws.Range("N15").Value = -5199.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8335.458000000001
$ws.Range("I32").Value = 5317.9316
$ws.Range("K32").Value = 5317.9316
$ws.Range("M32").Value = -5030.9316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12265.444
$ws.Range("I45").Value = 13548.625
$ws.Range("K45").Value = 13548.625
$ws.Range("M45").Value = -13171.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1559.3334
$ws.Range("I116").Value = 1506
$ws.Range("J116").Value = 1906
$ws.Range("K116").Value = 1506
$ws.Range("L116").Value = 1906
$ws.Range("M116").Value = 788
$ws.Range("N116").Value = -6494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1559.3334
$ws.Range("I3").Value = 1506
$ws.Range("J3").Value = 1906
$ws.Range("K3").Value = 1506
$ws.Range("L3").Value = 1906
$ws.Range("M3").Value = -1392
$ws.Range("N3").Value = -2134

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1512.8125
$ws.Range("I86").Value = 1457.0834
$ws.Range("J86").Value = 1680
$ws.Range("K86").Value = 1457.0834
$ws.Range("L86").Value = 1680
$ws.Range("M86").Value = -334.0834
$ws.Range("N86").Value = -3926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1512.8125
$ws.Range("I89").Value = 1457.0834
$ws.Range("J89").Value = 1680
$ws.Range("K89").Value = 7285.416999999999
$ws.Range("L89").Value = 8400
$ws.Range("M89").Value = -1669.416999999999
$ws.Range("N89").Value = -19632

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34376.84
$ws.Range("I31").Value = 38793.074
$ws.Range("J31").Value = 4567.25
$ws.Range("K31").Value = 38793.074
$ws.Range("L31").Value = 4567.25
$ws.Range("M31").Value = -38498.074
$ws.Range("N31").Value = -5157.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 34376.84
$ws.Range("I34").Value = 38793.074
$ws.Range("J34").Value = 4567.25
$ws.Range("K34").Value = 38793.074
$ws.Range("L34").Value = 4567.25
$ws.Range("M34").Value = -38591.074
$ws.Range("N34").Value = -4971.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 13918.75
$ws.Range("I60").Value = 12230.667
$ws.Range("J60").Value = 18983
$ws.Range("K60").Value = 12230.667
$ws.Range("L60").Value = 18983
$ws.Range("M60").Value = -11719.667
$ws.Range("N60").Value = -20005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2793.375
$ws.Range("I94").Value = 2516.5
$ws.Range("J94").Value = 2885.6667
$ws.Range("K94").Value = 2516.5
$ws.Range("L94").Value = 2885.6667
$ws.Range("M94").Value = -2065.5
$ws.Range("N94").Value = -3787.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2498.8462
$ws.Range("I99").Value = 2404.0908
$ws.Range("K99").Value = 2404.0908
$ws.Range("M99").Value = -906.0907999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2498.8462
$ws.Range("I126").Value = 2404.0908
$ws.Range("K126").Value = 7212.2724
$ws.Range("M126").Value = -4742.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2843.3333
$ws.Range("I140").Value = 2612
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 7836
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = -2656
$ws.Range("N140").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2886.2144
$ws.Range("I80").Value = 2845.3333
$ws.Range("J80").Value = 2959.8
$ws.Range("K80").Value = 2845.3333
$ws.Range("L80").Value = 2959.8
$ws.Range("M80").Value = -1847.3333
$ws.Range("N80").Value = -4955.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2886.2144
$ws.Range("I83").Value = 2845.3333
$ws.Range("J83").Value = 2959.8
$ws.Range("K83").Value = 14226.6665
$ws.Range("L83").Value = 14799
$ws.Range("M83").Value = -9234.666499999999
$ws.Range("N83").Value = -24783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10812.27
$ws.Range("I126").Value = 13689.368
$ws.Range("K126").Value = 41068.104
$ws.Range("M126").Value = -38598.104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3817.9644
$ws.Range("I61").Value = 3556.7727
$ws.Range("K61").Value = 3556.7727
$ws.Range("M61").Value = -3354.7727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1409.1666
$ws.Range("I93").Value = 1295.8422
$ws.Range("K93").Value = 1295.8422
$ws.Range("M93").Value = -47.84220000000005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3817.9644
$ws.Range("I113").Value = 3556.7727
$ws.Range("K113").Value = 3556.7727
$ws.Range("M113").Value = -1386.7727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 79377
$ws.Range("J41").Value = 79377
$ws.Range("L41").Value = 79377
$ws.Range("N41").Value = -80157

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 30000
$ws.Range("J56").Value = 30000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
